# se modifica la data para hacer smoke diez en pre prod
$wb = $excel.ActiveWorkbook

# --- DatosCuenta sheet (1st sheet) ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokPreProdDiez"
$wsCuenta.Range("B2").Value = "SmokeNamePreProdDiez"
$wsCuenta.Range("C2").Value = 27100118
$wsCuenta.Range("D2").Value = 120

# --- DatosHogar sheet (2nd sheet) ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 639

# --- DatosMotor sheet (3rd sheet) ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMP021"
$wsMotor.Range("B2").Value = "ABC12SSMP021"
$wsMotor.Range("C2").Value = "ZAZ123SSMP021"

# --- DatosAP sheet (4th sheet) ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200121

# Move the active tab / selection to DatosAP (it becomes the selected sheet)
$wsAP.Select()
$wsAP.Range("A3").Select()
